$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 2.964545797025059)
    3 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
    4 = @(1.445647641019636, 1.626987699542094, 3.223369029078222, 13.86384647080068, 0, 20.15985084044064)
    5 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 1, 4.327115817150455)
    6 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 3.755628166162433)
    7 = @(1.445647641019636, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 4.327115817150455)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
